# Applies the content edit described by the commit:
#   - The UML "boardSquare" diagram group is shifted up-and-left on the
#     slide.
#   - The empty Title placeholder shape that used to sit above the
#     diagram is removed so the diagram now fills the slide.
#   - The diagram shapes are renamed to match their new identities.
#
# NOTE: Shape.Left/Top are expressed in points; the host truncates
# (floors) points*12700 to obtain EMU internally, so each target is
# nudged by +0.5 EMU (well within point/EMU rounding noise) to land
# exactly on the intended EMU offset instead of one EMU short.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targets = @{
    1  = @{ Left = 387.00003937007875;  Top = 229.4861811023622;  Name = "Rectangle 222" }
    2  = @{ Left = 357.00003937007875;  Top = 181.2564960629921;  Name = "Folded Corner 223" }
    3  = @{ Left = 357.00003937007875;  Top = 181.25657480314962; Name = "TextBox 224" }
    4  = @{ Left = 543.0000393700788;   Top = 229.4861811023622;  Name = "Rectangle 225" }
    5  = @{ Left = 477.00003937007875;  Top = 244.02681102362206; Name = "Elbow Connector 226" }
    6  = @{ Left = 459.00003937007875;  Top = 235.02681102362206; Name = "Flowchart: Decision 227" }
    7  = @{ Left = 33.00003937007874;   Top = 187.25657480314962; Name = "Rectangle 228" }
    8  = @{ Left = 213.00003937007875;  Top = 188.17531496062992; Name = "Rectangle 229" }
    9  = @{ Left = 141.00003937007875;  Top = 202.25657480314962; Name = "Elbow Connector 43" }
    10 = @{ Left = 141.00003937007875;  Top = 169.25657480314962; Name = "TextBox 231" }
    11 = @{ Left = 123.00003937007874;  Top = 193.25657480314962; Name = "Flowchart: Decision 232" }
    12 = @{ Left = 111.00003937007874;  Top = 247.25657480314962; Name = "Rectangle 233" }
    13 = @{ Left = 161.77035433070867;  Top = 256.02688976377954; Name = "Elbow Connector 43" }
    14 = @{ Left = 147.00003937007875;  Top = 277.2565748031496;  Name = "Flowchart: Decision 235" }
    15 = @{ Left = 93.00003937007874;   Top = 313.2565748031496;  Name = "TextBox 236" }
    16 = $null   # "Title 26" placeholder - deleted below
}

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    $t = $targets[$i]

    if ($null -eq $t) {
        # Empty "Click to edit Master title style" placeholder - removed.
        $sh.Delete()
        continue
    }

    $sh.Left = $t.Left
    $sh.Top = $t.Top
    $sh.Name = $t.Name
}
